$wb = $excel.ActiveWorkbook

# --- Sheet "questions" ---
$ws1 = $wb.Worksheets.Item(1)

# B36: "Binary" -> "binary"
$ws1.Range("B36").Value = "binary"

# New row 41: a question about the telescope
$ws1.Range("A41").Value = "Who developed telescope?"
$ws1.Range("B41").Value = "Galileo Galilei"
$ws1.Range("C41").Value = $true
$ws1.Range("D41").Value = 2

# Make "questions" the active sheet/tab (was "archives questions")
$ws1.Activate()

# --- Sheet "students" ---
$ws3 = $wb.Worksheets.Item(3)

# Remove the placeholder/sample row 13 ("2xch1" / "sample"), shifting rows 14.. up by one
$ws3.Rows.Item(13).Delete()

# Append two new student rows at the end of the table
$ws3.Range("A14").Value = "09dqm"
$ws3.Range("B14").Value = "Sesgundo, Ryann Kim M"
$ws3.Range("C14").Value = "15e2b0d3c33891ebb0f1ef609ec419420c20e320ce94c65fbc8c3312448eb225"
$ws3.Range("D14").Value = 16

$ws3.Range("A15").Value = "p0rb2"
$ws3.Range("B15").Value = "Sesgundo, Ryann Kim"
$ws3.Range("C15").Value = "15e2b0d3c33891ebb0f1ef609ec419420c20e320ce94c65fbc8c3312448eb225"
$ws3.Range("D15").Value = 15
